$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.094.02'
$ws.Range("E2").Value = '  -1.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.421.22'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.99'
$ws.Range("E5").Value = '  -2.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.56'
$ws.Range("E6").Value = '  -2.26%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.528'
$ws.Range("E8").Value = '  -1.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.406.17'
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.08'
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.12'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("E15").Value = '  -2.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.835.36'
$ws.Range("E16").Value = '  -1.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.973.42'
$ws.Range("E17").Value = '  -1.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.393.48'
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.87'
$ws.Range("E19").Value = '  +8.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.63'
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.03'
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("E23").Value = '  +1.85%  '
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("E25").Value = '  -4.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.66'
$ws.Range("E26").Value = '  -1.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '585.72'
$ws.Range("E27").Value = '  -2.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.24'
$ws.Range("E28").Value = '  -9.36%  '
$ws.Range("E29").Value = '  -1.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0932'
$ws.Range("E30").Value = '  -3.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.88'
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("E32").Value = '  -5.25%  '
$ws.Range("E33").Value = '  -3.99%  '
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  -1.12%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.61'
$ws.Range("E37").Value = '  -5.52%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '151.57'
$ws.Range("E38").Value = '  -0.95%  '
$ws.Range("E39").Value = '  -2.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.24'
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.13'
$ws.Range("E44").Value = '  -4.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.35'
$ws.Range("E45").Value = '  -7.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0275'
$ws.Range("E46").Value = '  +3.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '143.17'
$ws.Range("E47").Value = '  +0.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.52'
$ws.Range("E48").Value = '  -2.92%  '
$ws.Range("E49").Value = '  -2.32%  '
$ws.Range("E50").Value = '  -2.07%  '
$ws.Range("E51").Value = '  -3.28%  '
